{"js": "// Replace the two-digit multiplication expressions in the table with the\n// updated set of problems/answers. Each old value is unique in the\n// document, so a targeted search + replace keeps the original run\n// formatting (font, size) intact for every cell.\nconst replacements = [\n  [\"43\u00d734=1462\", \"29\u00d782=2378\"],\n  [\"28\u00d750=1400\", \"18\u00d726=468\"],\n  [\"81\u00d714=1134\", \"47\u00d740=1880\"],\n  [\"74\u00d788=6512\", \"48\u00d798=4704\"],\n  [\"31\u00d785=2635\", \"97\u00d741=3977\"],\n  [\"29\u00d795=2755\", \"72\u00d740=2880\"],\n  [\"15\u00d734=510\", \"90\u00d772=6480\"],\n  [\"76\u00d797=7372\", \"67\u00d758=3886\"],\n  [\"55\u00d771=3905\", \"52\u00d721=1092\"],\n  [\"71\u00d733=2343\", \"90\u00d739=3510\"],\n  [\"97\u00d762=6014\", \"92\u00d783=7636\"],\n  [\"12\u00d755=660\", \"58\u00d714=812\"],\n  [\"52\u00d765=3380\", \"85\u00d718=1530\"],\n  [\"11\u00d783=913\", \"13\u00d764=832\"],\n  [\"80\u00d768=5440\", \"17\u00d757=969\"],\n  [\"72\u00d755=3960\", \"53\u00d767=3551\"],\n  [\"52\u00d793=4836\", \"36\u00d755=1980\"],\n  [\"53\u00d726=1378\", \"78\u00d785=6630\"],\n  [\"70\u00d716=1120\", \"45\u00d732=1440\"],\n  [\"11\u00d729=319\", \"59\u00d753=3127\"],\n  [\"76\u00d740=3040\", \"77\u00d781=6237\"],\n  [\"29\u00d768=1972\", \"86\u00d739=3354\"],\n  [\"33\u00d777=2541\", \"28\u00d721=588\"],\n  [\"80\u00d774=5920\", \"97\u00d799=9603\"],\n  [\"66\u00d741=2706\", \"51\u00d762=3162\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the two-digit multiplication expressions in the table with the\n# updated set of problems/answers. Each old value is unique in the\n# document, so Find/Replace (ReplaceAll) keeps the original run\n# formatting (font, size) intact for every cell.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"43\u00d734=1462\", \"29\u00d782=2378\"),\n    @(\"28\u00d750=1400\", \"18\u00d726=468\"),\n    @(\"81\u00d714=1134\", \"47\u00d740=1880\"),\n    @(\"74\u00d788=6512\", \"48\u00d798=4704\"),\n    @(\"31\u00d785=2635\", \"97\u00d741=3977\"),\n    @(\"29\u00d795=2755\", \"72\u00d740=2880\"),\n    @(\"15\u00d734=510\", \"90\u00d772=6480\"),\n    @(\"76\u00d797=7372\", \"67\u00d758=3886\"),\n    @(\"55\u00d771=3905\", \"52\u00d721=1092\"),\n    @(\"71\u00d733=2343\", \"90\u00d739=3510\"),\n    @(\"97\u00d762=6014\", \"92\u00d783=7636\"),\n    @(\"12\u00d755=660\", \"58\u00d714=812\"),\n    @(\"52\u00d765=3380\", \"85\u00d718=1530\"),\n    @(\"11\u00d783=913\", \"13\u00d764=832\"),\n    @(\"80\u00d768=5440\", \"17\u00d757=969\"),\n    @(\"72\u00d755=3960\", \"53\u00d767=3551\"),\n    @(\"52\u00d793=4836\", \"36\u00d755=1980\"),\n    @(\"53\u00d726=1378\", \"78\u00d785=6630\"),\n    @(\"70\u00d716=1120\", \"45\u00d732=1440\"),\n    @(\"11\u00d729=319\", \"59\u00d753=3127\"),\n    @(\"76\u00d740=3040\", \"77\u00d781=6237\"),\n    @(\"29\u00d768=1972\", \"86\u00d739=3354\"),\n    @(\"33\u00d777=2541\", \"28\u00d721=588\"),\n    @(\"80\u00d774=5920\", \"97\u00d799=9603\"),\n    @(\"66\u00d741=2706\", \"51\u00d762=3162\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
